# Update the UndoRedoSequenceDiagram: rename the "AddressBook" domain
# object family to "WishBook" throughout the sequence diagram on slide 1.
#
# Each target string is only part of a run inside its shape's text frame,
# so we address it with TextRange.Characters(start, length) rather than
# setting the whole TextRange.Text. That preserves the sibling
# runs/paragraphs (and their own formatting, e.g. "BookParser", "undo",
# "()", "resetData", "(", ")") exactly as the authored edit did, instead
# of collapsing the shape's text into a single run.

function Get-ShapeById($Slide, $Id) {
    for ($i = 1; $i -le $Slide.Shapes.Count; $i++) {
        $shp = $Slide.Shapes.Item($i)
        if ($shp.Id -eq $Id) {
            return $shp
        }
    }
    throw "No shape with Id=$Id on slide"
}

function Replace-Substring($Shape, $Old, $New) {
    $tr = $Shape.TextFrame.TextRange
    $full = $tr.Text
    $idx = $full.IndexOf($Old)
    if ($idx -lt 0) {
        throw "Substring '$Old' not found in shape '$($Shape.Name)' (text: '$full')"
    }
    $sub = $tr.Characters($idx + 1, $Old.Length)
    $sub.Text = $New
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape id=16 "Rectangle 62": ":Address" / "BookParser" -> ":Wish" / "BookParser"
Replace-Substring (Get-ShapeById $s 16) ":Address" ":Wish"

# Shape id=79 "TextBox 78": "undo" / "AddressBook" / "()" -> "undo" / "WishBook" / "()"
Replace-Substring (Get-ShapeById $s 79) "AddressBook" "WishBook"

# Shape id=84 "Rectangle 62": ":" / "VersionedAddressBook" -> ":" / "VersionedWishBook"
Replace-Substring (Get-ShapeById $s 84) "VersionedAddressBook" "VersionedWishBook"

# Shape id=88 "TextBox 87": "resetData" / "(" / "ReadOnlyAddressBook" / ")" -> "resetData" / "(" / "ReadOnlyWishBook" / ")"
Replace-Substring (Get-ShapeById $s 88) "ReadOnlyAddressBook" "ReadOnlyWishBook"
